# Auto-generated edit script applying the cryptos.xlsx price/volume update
# described by the commit "Updated cryptos list on Thu Sep 14 09:41:30 UTC 2023
# with GitHub Actions". Updates Price (D) and Volume(1h) (E) cell text for
# each affected row on the active sheet, preserving each cell's original
# Text data type (the source data are inline/shared strings, not numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.336.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.85%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.621.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.36%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "

# Row 6
$ws.Range("E6").Value = "  +0.11%  "

# Row 7
$ws.Range("E7").Value = "  +0.23%  "

# Row 8
$ws.Range("E8").Value = "  -0.15%  "

# Row 9
$ws.Range("E9").Value = "  +0.12%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.08%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.55%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.847.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.52%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.616.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.03%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.29%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.17%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.341.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.95%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.39%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0725"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.56%  "

# Row 19
$ws.Range("E19").Value = "  +0.02%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.88%  "

# Row 21
$ws.Range("E21").Value = "  -0.40%  "

# Row 22
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "

# Row 24
$ws.Range("E24").Value = "  -2.59%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.49%  "

# Row 26
$ws.Range("E26").Value = "  +0.18%  "

# Row 27
$ws.Range("E27").Value = "  -1.73%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.34%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.26%  "

# Row 30
$ws.Range("E30").Value = "  +8.23%  "

# Row 31
$ws.Range("E31").Value = "  +0.01%  "

# Row 32
$ws.Range("E32").Value = "  +0.93%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.08%  "

# Row 34
$ws.Range("E34").Value = "  -0.34%  "

# Row 35
$ws.Range("E35").Value = "  +2.24%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.157.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.15%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0163"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.07%  "

# Row 38
$ws.Range("E38").Value = "  +0.52%  "

# Row 40
$ws.Range("E40").Value = "  +0.17%  "

# Row 41
$ws.Range("E41").Value = "  +0.56%  "

# Row 42
$ws.Range("E42").Value = "  +3.96%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.782"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.33%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.759.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.45%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.27%  "

# Row 46
$ws.Range("E46").Value = "  +0.62%  "

# Row 47
$ws.Range("E47").Value = "  -0.94%  "

# Row 48
$ws.Range("E48").Value = "  +0.54%  "

# Row 49
$ws.Range("E49").Value = "  +0.56%  "

# Row 50
$ws.Range("E50").Value = "  -0.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.32"
$ws.Range("D51").Style = "Normal"

